$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item(1)
# Row 17
$ws.Range("H17").Value = 9092995
$ws.Range("J17").Value = 9092995
$ws.Range("L17").Value = 27278985
$ws.Range("N17").Value = -27279321
# Row 99
$ws.Range("H99").Value = 219.5
$ws.Range("I99").Value = 219.5
$ws.Range("K99").Value = 658.5
$ws.Range("M99").Value = 839.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item(2)
# Row 2
$ws.Range("H2").Value = 5834.1665
$ws.Range("I2").Value = 5285.5713
$ws.Range("K2").Value = 5285.5713
$ws.Range("M2").Value = -5172.5713
# Row 32
$ws.Range("H32").Value = 14327.611
$ws.Range("I32").Value = 12458.866
$ws.Range("J32").Value = 23671.334
$ws.Range("K32").Value = 12458.866
$ws.Range("L32").Value = 23671.334
$ws.Range("M32").Value = -12171.866
$ws.Range("N32").Value = -24245.334
# Row 61
$ws.Range("H61").Value = 8748.272000000001
$ws.Range("I61").Value = 8748.272000000001
$ws.Range("K61").Value = 8748.272000000001
$ws.Range("M61").Value = -8536.272000000001
# Row 74
$ws.Range("H74").Value = 4524.0835
$ws.Range("I74").Value = 2484.5715
$ws.Range("J74").Value = 11662.375
$ws.Range("K74").Value = 2484.5715
$ws.Range("L74").Value = 11662.375
$ws.Range("M74").Value = -1610.5715
$ws.Range("N74").Value = -13410.375
# Row 77
$ws.Range("H77").Value = 4524.0835
$ws.Range("I77").Value = 2484.5715
$ws.Range("J77").Value = 11662.375
$ws.Range("K77").Value = 12422.8575
$ws.Range("L77").Value = 58311.875
$ws.Range("M77").Value = -8054.8575
$ws.Range("N77").Value = -67047.875
# Row 116
$ws.Range("H116").Value = 5834.1665
$ws.Range("I116").Value = 5285.5713
$ws.Range("K116").Value = 5285.5713
$ws.Range("M116").Value = -2991.5713
# Row 122
$ws.Range("H122").Value = 1472.0555
$ws.Range("I122").Value = 1472.0555
$ws.Range("K122").Value = 4416.166499999999
$ws.Range("M122").Value = -1966.166499999999
# Row 132
$ws.Range("H132").Value = 1753.6389
$ws.Range("I132").Value = 1753.6389
$ws.Range("K132").Value = 5260.9167
$ws.Range("M132").Value = -2730.9167
# Row 136
$ws.Range("H136").Value = 8748.272000000001
$ws.Range("I136").Value = 8748.272000000001
$ws.Range("K136").Value = 26244.816
$ws.Range("M136").Value = -23694.816

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item(3)
# Row 3
$ws.Range("H3").Value = 5834.1665
$ws.Range("I3").Value = 5285.5713
$ws.Range("K3").Value = 5285.5713
$ws.Range("M3").Value = -5171.5713
# Row 86
$ws.Range("H86").Value = 1319.1818
$ws.Range("I86").Value = 1145.8125
$ws.Range("J86").Value = 1781.5
$ws.Range("K86").Value = 1145.8125
$ws.Range("L86").Value = 1781.5
$ws.Range("M86").Value = -22.8125
$ws.Range("N86").Value = -4027.5
# Row 89
$ws.Range("H89").Value = 1319.1818
$ws.Range("I89").Value = 1145.8125
$ws.Range("J89").Value = 1781.5
$ws.Range("K89").Value = 5729.0625
$ws.Range("L89").Value = 8907.5
$ws.Range("M89").Value = -113.0625
$ws.Range("N89").Value = -20139.5
# Row 99
$ws.Range("H99").Value = 4753.9585
$ws.Range("I99").Value = 3442.75
$ws.Range("K99").Value = 3442.75
$ws.Range("M99").Value = -1944.75
# Row 105
$ws.Range("H105").Value = 2210
$ws.Range("I105").Value = 2210
$ws.Range("K105").Value = 2210
$ws.Range("M105").Value = -463
# Row 141
$ws.Range("H141").Value = 45833.332
$ws.Range("J141").Value = 45833.332
$ws.Range("L141").Value = 45833.332
$ws.Range("N141").Value = -56193.332

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item(4)
# Row 16
$ws.Range("H16").Value = 3898.65
$ws.Range("I16").Value = 3273.4375
$ws.Range("K16").Value = 3273.4375
$ws.Range("M16").Value = -2986.4375
# Row 86
$ws.Range("H86").Value = 9141.200000000001
$ws.Range("I86").Value = 8591
$ws.Range("K86").Value = 8591
$ws.Range("M86").Value = -7468
# Row 89
$ws.Range("H89").Value = 9141.200000000001
$ws.Range("I89").Value = 8591
$ws.Range("K89").Value = 42955
$ws.Range("M89").Value = -37339
# Row 105
$ws.Range("H105").Value = 3131
$ws.Range("I105").Value = 3131
$ws.Range("K105").Value = 3131
$ws.Range("M105").Value = -1384
# Row 113
$ws.Range("H113").Value = 3898.65
$ws.Range("I113").Value = 3273.4375
$ws.Range("K113").Value = 3273.4375
$ws.Range("M113").Value = -1103.4375
# Row 140
$ws.Range("H140").Value = 93168.3
$ws.Range("J140").Value = 97886
$ws.Range("L140").Value = 97886
$ws.Range("N140").Value = -108246

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item(5)
# Row 2
$ws.Range("H2").Value = 84.75
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
# Row 86
$ws.Range("H86").Value = 2850.25
$ws.Range("J86").Value = 402
$ws.Range("L86").Value = 1206
$ws.Range("N86").Value = -3578
# Row 89
$ws.Range("H89").Value = 2850.25
$ws.Range("J89").Value = 402
$ws.Range("L89").Value = 3618
$ws.Range("N89").Value = -15474
# Row 112
$ws.Range("H112").Value = 5977
$ws.Range("I112").Value = 6846.3335
$ws.Range("J112").Value = 5325
$ws.Range("K112").Value = 20539.0005
$ws.Range("L112").Value = 15975
$ws.Range("M112").Value = -19431.0005
$ws.Range("N112").Value = -18191
# Row 120
$ws.Range("H120").Value = 16451
$ws.Range("I120").Value = 13811.8
$ws.Range("J120").Value = 19750
$ws.Range("K120").Value = 41435.39999999999
$ws.Range("L120").Value = 59250
$ws.Range("M120").Value = -36597.39999999999
$ws.Range("N120").Value = -68926
# Row 129
$ws.Range("H129").Value = 412791.97
$ws.Range("J129").Value = 533553.3
$ws.Range("L129").Value = 1600659.9
$ws.Range("N129").Value = -1610659.9

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item(6)
# Row 80
$ws.Range("H80").Value = 2161.2
$ws.Range("I80").Value = 2130.8333
$ws.Range("K80").Value = 2130.8333
$ws.Range("M80").Value = -1132.8333
# Row 83
$ws.Range("H83").Value = 2161.2
$ws.Range("I83").Value = 2130.8333
$ws.Range("K83").Value = 10654.1665
$ws.Range("M83").Value = -5662.166499999999
# Row 132
$ws.Range("H132").Value = 4180.3667
$ws.Range("I132").Value = 4225.0835
$ws.Range("K132").Value = 12675.2505
$ws.Range("M132").Value = -10145.2505

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item(7)
# Row 40
$ws.Range("H40").Value = 3291.32
$ws.Range("I40").Value = 3143.2778
$ws.Range("J40").Value = 3672
$ws.Range("K40").Value = 3143.2778
$ws.Range("L40").Value = 3672
$ws.Range("M40").Value = -3007.2778
$ws.Range("N40").Value = -3944
# Row 55
$ws.Range("H55").Value = 83333416
$ws.Range("I55").Value = 100000070
$ws.Range("J55").Value = 135
$ws.Range("K55").Value = 100000070
$ws.Range("L55").Value = 135
$ws.Range("M55").Value = -99999897
$ws.Range("N55").Value = -481
# Row 68
$ws.Range("H68").Value = 5067.364
$ws.Range("J68").Value = 8718.4
$ws.Range("L68").Value = 8718.4
$ws.Range("N68").Value = -10216.4
# Row 71
$ws.Range("H71").Value = 5067.364
$ws.Range("J71").Value = 8718.4
$ws.Range("L71").Value = 43592
$ws.Range("N71").Value = -51080
# Row 82
$ws.Range("H82").Value = 1026.8334
$ws.Range("J82").Value = 1082
$ws.Range("L82").Value = 1082
$ws.Range("N82").Value = -1804
# Row 85
$ws.Range("H85").Value = 1026.8334
$ws.Range("J85").Value = 1082
$ws.Range("L85").Value = 1082
$ws.Range("N85").Value = -3578

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item(8)
# Row 15
$ws.Range("H15").Value = 6991.6665
$ws.Range("J15").Value = 6989
$ws.Range("L15").Value = 6989
$ws.Range("N15").Value = -7565
# Row 122
$ws.Range("H122").Value = 2749.9666
$ws.Range("I122").Value = 1830
$ws.Range("K122").Value = 5490
$ws.Range("M122").Value = -3040

# ---- Deletions ----
$ws = $wb.Worksheets.Item(5)
$ws.Range("N2").ClearContents()
